# "Ran Priors with _1"
#
# Re-ran the Bayesian priors table with a wider prior variance / shrinkage
# setting: column D (prior variance) goes from 1 -> 5 and column E (prior
# mean shrink) goes from 0.15/0.05 -> 0.2 for every loading row (2-44).
# Rows 45-67 (the "A1" block) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D44").Value = 5
$ws.Range("E2:E44").Value = 0.2

# Match the page setup recorded for this run (portrait, letter/A4-class
# paper) and leave the selection on the column that was just edited.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("E2:E44").Select()
